# Apply updated Universalis market-price derived values to the Leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 32899.902
$ws.Range("J28").Value = 712
$ws.Range("L28").Value = 712
$ws.Range("N28").Value = -1682
$ws.Range("H43").Value = 824
$ws.Range("J43").Value = 800
$ws.Range("L43").Value = 800
$ws.Range("N43").Value = -938
$ws.Range("H62").Value = 36612.75
$ws.Range("J62").Value = 43984.617
$ws.Range("L62").Value = 43984.617
$ws.Range("N62").Value = -45232.617
$ws.Range("H65").Value = 36612.75
$ws.Range("J65").Value = 43984.617
$ws.Range("L65").Value = 219923.085
$ws.Range("N65").Value = -226163.085
$ws.Range("H96").Value = 185.09091
$ws.Range("I96").Value = 146.75
$ws.Range("J96").Value = 287.33334
$ws.Range("K96").Value = 440.25
$ws.Range("L96").Value = 862.0000200000001
$ws.Range("M96").Value = 932.75
$ws.Range("N96").Value = -3608.00002
$ws.Range("H98").Value = 3113.3547
$ws.Range("I98").Value = 3003.44
$ws.Range("J98").Value = 3571.3333
$ws.Range("K98").Value = 3003.44
$ws.Range("L98").Value = 3571.3333
$ws.Range("M98").Value = -1505.44
$ws.Range("N98").Value = -6567.3333
$ws.Range("H100").Value = 4471.35
$ws.Range("I100").Value = 5235.5835
$ws.Range("J100").Value = 3325
$ws.Range("K100").Value = 5235.5835
$ws.Range("L100").Value = 3325
$ws.Range("M100").Value = -4694.5835
$ws.Range("N100").Value = -4407
$ws.Range("H107").Value = 300.15384
$ws.Range("I107").Value = 306.5
$ws.Range("K107").Value = 306.5
$ws.Range("M107").Value = 1613.5
$ws.Range("H116").Value = 22483710
$ws.Range("I116").Value = 15699153
$ws.Range("J116").Value = 33339004
$ws.Range("K116").Value = 15699153
$ws.Range("L116").Value = 33339004
$ws.Range("M116").Value = -15695711
$ws.Range("N116").Value = -33345888
$ws.Range("H122").Value = 3113.3547
$ws.Range("I122").Value = 3003.44
$ws.Range("J122").Value = 3571.3333
$ws.Range("K122").Value = 9010.32
$ws.Range("L122").Value = 10713.9999
$ws.Range("M122").Value = -6560.32
$ws.Range("N122").Value = -15613.9999
$ws.Range("H132").Value = 7587.6113
$ws.Range("I132").Value = 8042.7144
$ws.Range("J132").Value = 5994.75
$ws.Range("K132").Value = 24128.1432
$ws.Range("L132").Value = 17984.25
$ws.Range("M132").Value = -21598.1432
$ws.Range("N132").Value = -23044.25
$ws.Range("H137").Value = 2881.7778
$ws.Range("I137").Value = 2311.913
$ws.Range("J137").Value = 3890
$ws.Range("K137").Value = 6935.739
$ws.Range("L137").Value = 11670
$ws.Range("M137").Value = -4385.739
$ws.Range("N137").Value = -16770
$ws.Range("H138").Value = 2223.3333
$ws.Range("I138").Value = 1997.2222
$ws.Range("J138").Value = 2392.9167
$ws.Range("K138").Value = 5991.6666
$ws.Range("L138").Value = 7178.750100000001
$ws.Range("M138").Value = -851.6665999999996
$ws.Range("N138").Value = -17458.7501
$ws.Range("H141").Value = 7268.9
$ws.Range("I141").Value = 7268.9
$ws.Range("K141").Value = 21806.7
$ws.Range("M141").Value = -16626.7

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H31").Value = 7500
$ws.Range("I31").Value = 5000
$ws.Range("J31").Value = 10000
$ws.Range("K31").Value = 5000
$ws.Range("L31").Value = 10000
$ws.Range("M31").Value = -4706
$ws.Range("N31").Value = -10588
$ws.Range("H32").Value = 5774.81
$ws.Range("I32").Value = 4909.743
$ws.Range("J32").Value = 18577.8
$ws.Range("K32").Value = 4909.743
$ws.Range("L32").Value = 18577.8
$ws.Range("M32").Value = -4622.743
$ws.Range("N32").Value = -19151.8
$ws.Range("H74").Value = 3812.7368
$ws.Range("I74").Value = 3264.5
$ws.Range("J74").Value = 5868.625
$ws.Range("K74").Value = 3264.5
$ws.Range("L74").Value = 5868.625
$ws.Range("M74").Value = -2390.5
$ws.Range("N74").Value = -7616.625
$ws.Range("H77").Value = 3812.7368
$ws.Range("I77").Value = 3264.5
$ws.Range("J77").Value = 5868.625
$ws.Range("K77").Value = 16322.5
$ws.Range("L77").Value = 29343.125
$ws.Range("M77").Value = -11954.5
$ws.Range("N77").Value = -38079.125
$ws.Range("H132").Value = 4715.615
$ws.Range("I132").Value = 3201.5715
$ws.Range("J132").Value = 6482
$ws.Range("K132").Value = 9604.7145
$ws.Range("L132").Value = 19446
$ws.Range("M132").Value = -7074.7145
$ws.Range("N132").Value = -24506

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 15690500
$ws.Range("I107").Value = 78303.69500000001
$ws.Range("K107").Value = 78303.69500000001
$ws.Range("M107").Value = -76383.69500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 57764.668
$ws.Range("I62").Value = 2206.6
$ws.Range("J62").Value = 127212.25
$ws.Range("K62").Value = 2206.6
$ws.Range("L62").Value = 127212.25
$ws.Range("M62").Value = -1582.6
$ws.Range("N62").Value = -128460.25
$ws.Range("H65").Value = 57764.668
$ws.Range("I65").Value = 2206.6
$ws.Range("J65").Value = 127212.25
$ws.Range("K65").Value = 11033
$ws.Range("L65").Value = 636061.25
$ws.Range("M65").Value = -7913
$ws.Range("N65").Value = -642301.25
$ws.Range("H99").Value = 5865.923
$ws.Range("I99").Value = 7163.5713
$ws.Range("K99").Value = 7163.5713
$ws.Range("M99").Value = -5665.5713
$ws.Range("H105").Value = 1094
$ws.Range("I105").Value = 1164.75
$ws.Range("K105").Value = 1164.75
$ws.Range("M105").Value = 582.25
$ws.Range("H107").Value = 493.87097
$ws.Range("J107").Value = 620.5
$ws.Range("L107").Value = 620.5
$ws.Range("N107").Value = -4460.5
$ws.Range("H126").Value = 5865.923
$ws.Range("I126").Value = 7163.5713
$ws.Range("K126").Value = 21490.7139
$ws.Range("M126").Value = -19020.7139
$ws.Range("H134").Value = 4593.5356
$ws.Range("I134").Value = 4129.6665
$ws.Range("K134").Value = 12388.9995
$ws.Range("M134").Value = -9853.999500000002

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 35707304
$ws.Range("I4").Value = 35707304
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 107121912
$ws.Range("L4").Value = 0
$ws.Range("M4").ClearContents()
$ws.Range("N4").Value = -107121800
$ws.Range("H11").Value = 116511.695
$ws.Range("I11").Value = 31336.5
$ws.Range("K11").Value = 94009.5
$ws.Range("M11").Value = -93869.5
$ws.Range("H81").Value = 7785.273
$ws.Range("I81").Value = 7671
$ws.Range("J81").Value = 7828.125
$ws.Range("K81").Value = 23013
$ws.Range("L81").Value = 23484.375
$ws.Range("M81").Value = -21890
$ws.Range("N81").Value = -25730.375
$ws.Range("H84").Value = 7785.273
$ws.Range("I84").Value = 7671
$ws.Range("J84").Value = 7828.125
$ws.Range("K84").Value = 69039
$ws.Range("L84").Value = 70453.125
$ws.Range("M84").Value = -63423
$ws.Range("N84").Value = -81685.125
$ws.Range("H122").Value = 618.1667
$ws.Range("I122").Value = 613.5714
$ws.Range("K122").Value = 5522.1426
$ws.Range("M122").Value = -3072.1426
$ws.Range("H129").Value = 2308.3333
$ws.Range("J129").Value = 2707.5
$ws.Range("L129").Value = 8122.5
$ws.Range("N129").Value = -18122.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 10470.125
$ws.Range("I70").Value = 9753
$ws.Range("K70").Value = 9753
$ws.Range("M70").Value = -9483
$ws.Range("H73").Value = 10470.125
$ws.Range("I73").Value = 9753
$ws.Range("K73").Value = 9753
$ws.Range("M73").Value = -8817
$ws.Range("H99").Value = 11354
$ws.Range("I99").Value = 10942.5
$ws.Range("K99").Value = 10942.5
$ws.Range("M99").Value = -8696.5
$ws.Range("H102").Value = 2473.9033
$ws.Range("I102").Value = 1561.2106
$ws.Range("K102").Value = 1561.2106
$ws.Range("M102").Value = 60.78939999999989
$ws.Range("H117").Value = 39154.5
$ws.Range("J117").Value = 39154.5
$ws.Range("L117").Value = 39154.5
$ws.Range("N117").Value = -46038.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 2422.75
$ws.Range("I46").Value = 825.5
$ws.Range("K46").Value = 825.5
$ws.Range("M46").Value = -637.5
$ws.Range("H93").Value = 1851.0454
$ws.Range("I93").Value = 1807.3334
$ws.Range("K93").Value = 1807.3334
$ws.Range("M93").Value = -559.3334
$ws.Range("H122").Value = 3942.0833
$ws.Range("I122").Value = 3571.5881
$ws.Range("K122").Value = 10714.7643
$ws.Range("M122").Value = -8264.764299999999

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29999
$ws.Range("I64").Value = 29999
$ws.Range("K64").Value = 29999
$ws.Range("M64").Value = -29751
$ws.Range("H67").Value = 29999
$ws.Range("I67").Value = 29999
$ws.Range("K67").Value = 29999
$ws.Range("M67").Value = -29141
$ws.Range("H96").Value = 2342.889
$ws.Range("J96").Value = 2665.6667
$ws.Range("L96").Value = 2665.6667
$ws.Range("N96").Value = -5411.6667
$ws.Range("H113").Value = 1339.9333
$ws.Range("I113").Value = 443.88235
$ws.Range("J113").Value = 2511.6924
$ws.Range("K113").Value = 1331.64705
$ws.Range("L113").Value = 7535.0772
$ws.Range("M113").Value = 838.35295
$ws.Range("N113").Value = -11875.0772
$ws.Range("H132").Value = 492093.8
$ws.Range("I132").Value = 517081.94
$ws.Range("J132").Value = 4825
$ws.Range("K132").Value = 1551245.82
$ws.Range("L132").Value = 14475
$ws.Range("M132").Value = -1548715.82
$ws.Range("N132").Value = -19535
$ws.Range("H136").Value = 5485.968
$ws.Range("I136").Value = 5994.5
$ws.Range("K136").Value = 17983.5
$ws.Range("M136").Value = -15433.5
$ws.Range("H140").Value = 115999.664
$ws.Range("I140").Value = 68999
$ws.Range("J140").Value = 139500
$ws.Range("K140").Value = 68999
$ws.Range("L140").Value = 139500
$ws.Range("M140").Value = -63819
